$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("H8").Value = 467.14285
$ws.Range("I8").Value = 533.6667
$ws.Range("K8").Value = 1601.0001
$ws.Range("M8").Value = -1462.0001
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("H15").Value = 1787.6415
$ws.Range("I15").Value = 1787.6415
$ws.Range("K15").Value = 5362.9245
$ws.Range("M15").Value = -5193.9245
$ws.Range("H51").Value = 6500
$ws.Range("I51").Value = 6500
$ws.Range("K51").Value = 6500
$ws.Range("M51").Value = -6016
$ws.Range("H55").Value = 319.64706
$ws.Range("J55").Value = 411.625
$ws.Range("L55").Value = 411.625
$ws.Range("N55").Value = -839.625
$ws.Range("H57").Value = 59885.75
$ws.Range("J57").Value = 59885.75
$ws.Range("L57").Value = 179657.25
$ws.Range("N57").Value = -180655.25
$ws.Range("H62").Value = 13560.071
$ws.Range("I62").Value = 11542.9375
$ws.Range("K62").Value = 11542.9375
$ws.Range("M62").Value = -10918.9375
$ws.Range("H65").Value = 13560.071
$ws.Range("I65").Value = 11542.9375
$ws.Range("K65").Value = 57714.6875
$ws.Range("M65").Value = -54594.6875
$ws.Range("H113").Value = 3200
$ws.Range("I113").Value = 2933.3333
$ws.Range("J113").Value = 3466.6667
$ws.Range("K113").Value = 2933.3333
$ws.Range("L113").Value = 3466.6667
$ws.Range("M113").Value = 320.6667000000002
$ws.Range("N113").Value = -9974.6667
$ws.Range("H132").Value = 61975.49
$ws.Range("I132").Value = 70236.64999999999
$ws.Range("K132").Value = 210709.95
$ws.Range("M132").Value = -208179.95
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15028.6
$ws.Range("I2").Value = 6241
$ws.Range("K2").Value = 6241
$ws.Range("M2").Value = -6128
$ws.Range("H74").Value = 3731.9285
$ws.Range("I74").Value = 916.3333
$ws.Range("J74").Value = 8800
$ws.Range("K74").Value = 916.3333
$ws.Range("L74").Value = 8800
$ws.Range("M74").Value = -42.33330000000001
$ws.Range("N74").Value = -10548
$ws.Range("H77").Value = 3731.9285
$ws.Range("I77").Value = 916.3333
$ws.Range("J77").Value = 8800
$ws.Range("K77").Value = 4581.6665
$ws.Range("L77").Value = 44000
$ws.Range("M77").Value = -213.6665000000003
$ws.Range("N77").Value = -52736
$ws.Range("H116").Value = 15028.6
$ws.Range("I116").Value = 6241
$ws.Range("K116").Value = 6241
$ws.Range("M116").Value = -3947
$ws.Range("H122").Value = 3030.75
$ws.Range("I122").Value = 2844.7144
$ws.Range("K122").Value = 8534.143199999999
$ws.Range("M122").Value = -6084.143199999999
$ws.Range("H132").Value = 590050.5
$ws.Range("I132").Value = 691525.0600000001
$ws.Range("K132").Value = 2074575.18
$ws.Range("M132").Value = -2072045.18
$ws.Range("H135").Value = 69846
$ws.Range("J135").Value = 69846
$ws.Range("L135").Value = 69846
$ws.Range("N135").Value = -79986

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15028.6
$ws.Range("I3").Value = 6241
$ws.Range("K3").Value = 6241
$ws.Range("M3").Value = -6127
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("H116").Value = 56348.4
$ws.Range("J116").Value = 56348.4
$ws.Range("L116").Value = 56348.4
$ws.Range("N116").Value = -65526.4
$ws.Range("N60").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18776.482
$ws.Range("I31").Value = 6558.6665
$ws.Range("J31").Value = 43212.11
$ws.Range("K31").Value = 6558.6665
$ws.Range("L31").Value = 43212.11
$ws.Range("M31").Value = -6263.6665
$ws.Range("N31").Value = -43802.11
$ws.Range("H34").Value = 18776.482
$ws.Range("I34").Value = 6558.6665
$ws.Range("J34").Value = 43212.11
$ws.Range("K34").Value = 6558.6665
$ws.Range("L34").Value = 43212.11
$ws.Range("M34").Value = -6356.6665
$ws.Range("N34").Value = -43616.11
$ws.Range("H41").Value = 15853.385
$ws.Range("I41").Value = 3039.3333
$ws.Range("J41").Value = 19697.6
$ws.Range("K41").Value = 3039.3333
$ws.Range("L41").Value = 19697.6
$ws.Range("M41").Value = -2611.3333
$ws.Range("N41").Value = -20553.6
$ws.Range("H97").Value = 52500
$ws.Range("J97").Value = 52500
$ws.Range("L97").Value = 52500
$ws.Range("N97").Value = -54482
$ws.Range("H134").Value = 2420.5417
$ws.Range("I134").Value = 2511.8
$ws.Range("J134").Value = 1964.25
$ws.Range("K134").Value = 7535.400000000001
$ws.Range("L134").Value = 5892.75
$ws.Range("M134").Value = -5000.400000000001
$ws.Range("N134").Value = -10962.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H132").Value = 546.7143
$ws.Range("I132").Value = 526.1667
$ws.Range("K132").Value = 4735.5003
$ws.Range("M132").Value = -2205.5003
$ws.Range("H136").Value = 6490.3335
$ws.Range("I136").Value = 4697.7144
$ws.Range("K136").Value = 14093.1432
$ws.Range("M136").Value = -8993.143199999999
$ws.Range("H138").Value = 4182.375
$ws.Range("I138").Value = 4641.8
$ws.Range("J138").Value = 3416.6667
$ws.Range("K138").Value = 13925.4
$ws.Range("L138").Value = 10250.0001
$ws.Range("M138").Value = -8785.400000000001
$ws.Range("N138").Value = -20530.0001
$ws.Range("M106").ClearContents()
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11157178
$ws.Range("I11").Value = 16679117
$ws.Range("J11").Value = 113300
$ws.Range("K11").Value = 16679117
$ws.Range("L11").Value = 113300
$ws.Range("M11").Value = -16678978
$ws.Range("N11").Value = -113578
$ws.Range("H70").Value = 4876.222
$ws.Range("I70").Value = 5063.3335
$ws.Range("K70").Value = 5063.3335
$ws.Range("M70").Value = -4793.3335
$ws.Range("H73").Value = 4876.222
$ws.Range("I73").Value = 5063.3335
$ws.Range("K73").Value = 5063.3335
$ws.Range("M73").Value = -4127.3335
$ws.Range("H122").Value = 4593.268
$ws.Range("J122").Value = 4995.852
$ws.Range("L122").Value = 14987.556
$ws.Range("N122").Value = -19887.556

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 145.5
$ws.Range("I9").Value = 145.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 145.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 78.5
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("H46").Value = 6438.769
$ws.Range("I46").Value = 9463.375
$ws.Range("K46").Value = 9463.375
$ws.Range("M46").Value = -9275.375
$ws.Range("N9").ClearContents()
$ws.Range("M13").ClearContents()
